$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "44.589.88"
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = "  +4.33%  "

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.287.08"
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = "  +3.11%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.22%  "

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "321.25"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +1.60%  "

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "107.07"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +7.52%  "

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.594"
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = "  +0.44%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.21%  "

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.576"
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +2.41%  "

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "39.03"
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +5.21%  "

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0848"
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = "  +2.55%  "

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.96"
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +2.35%  "

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.108"
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  +1.49%  "

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.891"
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +3.30%  "

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.629.99"
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = "  +2.70%  "

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "14.71"
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = "  +3.20%  "

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.279.80"
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  +2.62%  "

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "44.399.98"
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  +3.85%  "

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "14.23"
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -5.23%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +5.13%  "

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.57"
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +2.17%  "

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "66.70"
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  +2.03%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +2.35%  "

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "240.22"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +2.00%  "

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.22"
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +3.81%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.06%  "

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.28"
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +1.92%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "Toncoin"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  +0.80%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "38.57"
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +13.13%  "

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.60"
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +3.53%  "

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.80"
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +1.73%  "

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "163.45"
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  +4.52%  "

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0891"
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -0.50%  "

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.75"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -0.90%  "

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.08"
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +5.75%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +2.09%  "

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.118"
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +13.86%  "

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.122"
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -0.21%  "

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.97"
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +2.97%  "

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.50"
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +1.57%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "Celestia"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.66"
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  +26.99%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0330"
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +1.69%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.02%  "

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.782.27"
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -7.73%  "

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.210"
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +1.40%  "

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "87.12"
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -1.71%  "

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.51"
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +2.49%  "

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "60.85"
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  +0.74%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Stacks"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.73"
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +7.13%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "ordi"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "75.35"
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = "  +0.32%  "

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "105.13"
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +2.50%  "
